# Update for first draft
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# New header cells (F1:H1) need the same bold/bordered/centered style as
# the existing header cells (B1:E1) - copy format from E1 across.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows ---
$algos = @("LR", "LDA", "KNN", "DTREE", "RTREE", "XTREE", "SVM")

$data = @(
    @(0.7581038730459106, 0.02621230083945392, 0.5612070710708444, 0.03503521274125497, 0.7257261574406748, 0.01818538711038182),
    @(0.7538638742612394, 0.0288019356229122,  0.5398090777497146, 0.03353151315474634, 0.726662142971398,  0.02246338337303938),
    @(0.7812755034078485, 0.02120015955865973, 0.6672942846635481, 0.02809822106611554, 0.74285078278352,   0.01240178660086847),
    @(0.7109832275150931, 0.02098716405844013, 0.5908809514844836, 0.03353551269733934, 0.6805852730898864, 0.02543286732970798),
    @(0.7006494510868156, 0.02971012758622779, 0.5161305134198356, 0.02505603624963511, 0.6598906605537163, 0.01827512639147544),
    @(0.7769181335858562, 0.02952922178989007, 0.612156124748743,  0.03046591651716983, 0.747732333759792,  0.02132693641557746),
    @(0.7869941765067398, 0.0293372834328052,  0.6595700072994649, 0.0348183318750939,  0.7540670400568846, 0.02015347549556609)
)

for ($i = 0; $i -lt $algos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $algos[$i]
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($row, 3 + $j).Value = $data[$i][$j]
    }
}

# Remove the now-unused 8th data row (row 9) left over from the old 8-row table
$ws.Rows.Item(9).Delete()
